$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 9: entrada de mercancia / inventory record updated ---
$ws.Range("B9").Value = "399840300"
$ws.Range("C9").Value = "3052754289"
$ws.Range("D9").Value = "732111324709674"

# --- Row 10 ---
$ws.Range("B10").Value = "592168140"
$ws.Range("C10").Value = "3046010569"
$ws.Range("D10").Value = "732111193280551"

# --- Row 11 ---
$ws.Range("B11").Value = "399840300"
$ws.Range("C11").Value = "3052754289"
$ws.Range("D11").Value = "732111324709674"

# --- Row 12 ---
$ws.Range("B12").Value = "313778543"
$ws.Range("C12").Value = "3052754293"
$ws.Range("D12").Value = "732111324709675"
$ws.Range("E12").Value = "cliente nit a nit"

# --- Row 13 ---
$ws.Range("B13").Value = "994114856"
$ws.Range("C13").Value = "3052754321"
$ws.Range("D13").Value = "732111324709676"
$ws.Range("E13").Value = "901963208"

# --- Row 14: values change and stray "horizontal=left" style is dropped (style 7 -> style 1) ---
$ws.Range("B14:C14").ClearFormats()
$ws.Range("B14:C14").NumberFormat = "@"
$ws.Range("B14").Value = "721106626"
$ws.Range("C14").Value = "3046008586"
$ws.Range("D14").Value = "732111193278871"

# --- New rows 15 & 16, same button layout as the rows above them ---
$ws.Range("A15:D16").NumberFormat = "@"

$ws.Range("A15").Value = "10960370"
$ws.Range("B15").Value = "61962571"
$ws.Range("C15").Value = "3046008587"
$ws.Range("D15").Value = "732111193280535"

$ws.Range("A16").Value = "10960370"
$ws.Range("B16").Value = "987388666"
$ws.Range("C16").Value = "3045987650"
$ws.Range("D16").Value = "732111193278858"

# --- Restore the selection Excel leaves behind on this sheet ---
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
